$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data lives in columns B:M with an unused column A in front
# of it (dimension starts at B1). Deleting column A shifts the whole
# table one column to the left (B->A, C->B, ..., M->L), which is what
# the diff shows: headers/values move left by one column, the
# dimension becomes A1:L4, and the bestFit column-width formatting
# moves from column B to column A.
$ws.Columns.Item(1).Delete() | Out-Null

# Matches the saved selection recorded in the target file.
$ws.Range("K17").Select() | Out-Null
